# edit.ps1
# Applies the OOXML diff to the document using Word COM interop semantics.
#
# Summary of changes:
#  1. Split the "Summarize" / "What Marketing Teams Achieve..." paragraph
#     into two separate paragraphs; the "Summarize" run loses its
#     Heading1Char rStyle, and the second paragraph's run gets a
#     lightGray highlight.
#  2. Add a <w:lastRenderedPageBreak/> before the dashed divider line run.
#  3. Split the "Your sales team's feedback..." run into two runs (break
#     inserted mid-sentence before "know which content...") with a
#     <w:lastRenderedPageBreak/> on the second run.
#  4. Remove the <w:lastRenderedPageBreak/> that used to precede
#     "Accelerate Campaign Planning and Execution".

$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParagraphRangeForText($searchText) {
    $range = $d.Content
    $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $para = $range.Paragraphs(1)
    return $para.Range
}

# ---------------------------------------------------------------------
# Change 1: "Summarize" heading paragraph split into two paragraphs.
# ---------------------------------------------------------------------
$pr1 = Get-ParagraphRangeForText("Summarize")
$xml1 = "<w:p $W>" +
          "<w:pPr>" +
            "<w:pStyle w:val='Heading3'/>" +
            "<w:spacing w:before='281' w:after='281'/>" +
          "</w:pPr>" +
          "<w:r><w:t>Summarize</w:t></w:r>" +
        "</w:p>" +
        "<w:p $W>" +
          "<w:pPr>" +
            "<w:pStyle w:val='Heading3'/>" +
            "<w:spacing w:before='281' w:after='281'/>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii='Aptos' w:eastAsia='Aptos' w:hAnsi='Aptos' w:cs='Aptos'/>" +
              "<w:sz w:val='24'/>" +
              "<w:szCs w:val='24'/>" +
            "</w:rPr>" +
          "</w:pPr>" +
          "<w:r><w:br/></w:r>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii='Arial' w:eastAsia='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
              "<w:b/>" +
              "<w:bCs/>" +
              "<w:color w:val='000000' w:themeColor='text1'/>" +
              "<w:highlight w:val='lightGray'/>" +
            "</w:rPr>" +
            "<w:t>What Marketing Teams Achieve with Synoptix AI (Callout box)</w:t>" +
          "</w:r>" +
        "</w:p>"
$pr1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Change 2: Dashed divider line gains a lastRenderedPageBreak.
# ---------------------------------------------------------------------
$dashes = "----------------------------------------------------------------------------------------------------------------"
$pr2 = Get-ParagraphRangeForText($dashes)
$xml2 = "<w:p $W>" +
          "<w:pPr>" +
            "<w:spacing w:before='240' w:after='240'/>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii='Arial' w:eastAsia='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
              "<w:color w:val='000000' w:themeColor='text1'/>" +
            "</w:rPr>" +
          "</w:pPr>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii='Arial' w:eastAsia='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
              "<w:color w:val='000000' w:themeColor='text1'/>" +
            "</w:rPr>" +
            "<w:lastRenderedPageBreak/>" +
            "<w:t>$dashes</w:t>" +
          "</w:r>" +
        "</w:p>"
$pr2.InsertXML($xml2)

# ---------------------------------------------------------------------
# Change 3: Split the "Your sales team's feedback..." run into two runs
# with a lastRenderedPageBreak on the second.
# ---------------------------------------------------------------------
$pr3 = Get-ParagraphRangeForText("Your sales team")
$part1 = [char]0x2019 + "s feedback is a goldmine, and Synoptix ensures your marketing team can tap into it. With integrations across CRM, call platforms, and sales tools, Synoptix helps you see what" + [char]0x2019 + "s working, what" + [char]0x2019 + "s missing, and what needs improvement. You" + [char]0x2019 + "ll "
$textRun1 = "Your sales team" + $part1
$textRun2 = "know which content gets used, what questions buyers are asking, and how to adapt campaigns without waiting "
$xml3 = "<w:p $W>" +
          "<w:pPr>" +
            "<w:spacing w:before='240' w:after='240'/>" +
          "</w:pPr>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii='Arial' w:eastAsia='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
              "<w:color w:val='000000' w:themeColor='text1'/>" +
            "</w:rPr>" +
            "<w:t xml:space='preserve'>$textRun1</w:t>" +
          "</w:r>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii='Arial' w:eastAsia='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
              "<w:color w:val='000000' w:themeColor='text1'/>" +
            "</w:rPr>" +
            "<w:lastRenderedPageBreak/>" +
            "<w:t xml:space='preserve'>$textRun2</w:t>" +
          "</w:r>" +
          "<w:proofErr w:type='gramStart'/>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii='Arial' w:eastAsia='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
              "<w:color w:val='000000' w:themeColor='text1'/>" +
            "</w:rPr>" +
            "<w:t>on</w:t>" +
          "</w:r>" +
          "<w:proofErr w:type='gramEnd'/>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii='Arial' w:eastAsia='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
              "<w:color w:val='000000' w:themeColor='text1'/>" +
            "</w:rPr>" +
            "<w:t xml:space='preserve'> a quarterly review.</w:t>" +
          "</w:r>" +
        "</w:p>"
$pr3.InsertXML($xml3)

# ---------------------------------------------------------------------
# Change 4: Remove lastRenderedPageBreak before "Accelerate Campaign
# Planning and Execution".
# ---------------------------------------------------------------------
$pr4 = Get-ParagraphRangeForText("Accelerate Campaign Planning and Execution")
$xml4 = "<w:p $W>" +
          "<w:pPr>" +
            "<w:pStyle w:val='Heading2'/>" +
            "<w:spacing w:before='360'/>" +
          "</w:pPr>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii='Arial' w:eastAsia='Arial' w:hAnsi='Arial' w:cs='Arial'/>" +
              "<w:b/>" +
              "<w:bCs/>" +
              "<w:color w:val='000000' w:themeColor='text1'/>" +
              "<w:sz w:val='34'/>" +
              "<w:szCs w:val='34'/>" +
            "</w:rPr>" +
            "<w:t>Accelerate Campaign Planning and Execution</w:t>" +
          "</w:r>" +
        "</w:p>"
$pr4.InsertXML($xml4)

Write-Host "Edits applied successfully."
